$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: populate cells that introduce brand-new shared strings, in the exact
#     order needed so the saved sharedStrings table lines up with the target indices ---
$ws.Range('A3').Value = 'Commuter Challenge'
$ws.Range('B2').Value = 'factitious'
$ws.Range('B3').Value = 'commuter'
$ws.Range('E1').Value = 'client'
$ws.Range('F1').Value = 'myRoles'
$ws.Range('G1').Value = 'brief'
$ws.Range('H1').Value = 'product'
$ws.Range('I1').Value = 'outcome'
$ws.Range('E2').Value = 'AU Game Lab'
$ws.Range('F2').Value = 'Game design, UX design, UI design'
$ws.Range('G2').Value = 'In the wake of the 2016 election, there was a huge discussion in the journalism community about how readers evaluate truth and trustworthiness in a news outlet. The client was interested in conducting a research project about readers’ ability to identify “truthiness” in a news story, and wanted to get as many plays as possible.'
$ws.Range('H2').Value = 'Using game design principles, we designed an interface that made news story evaluation fun. Users were presented with a new story, and a Tinder-like swipe mechanic let them designate a story as “real” or “fake” news.'
$ws.Range('I2').Value = 'Factitious was played over 1.6 million times in the first three days it was released and was covered in several major media outlets. The client was able to collect enough data for their own research, as well as ultimately making the core software available as an open source tool for playful polling systems.'
$ws.Range('A4').Value = 'Comics MFA Digital Anthology'
$ws.Range('B4').Value = 'cca'
$ws.Range('C2').Value = 'Newsgame'
$ws.Range('E4').Value = 'California College of the Arts'
$ws.Range('F4').Value = 'Web app design'
$ws.Range('C4').Value = 'Web app'
$ws.Range('E3').Value = 'asdf'

# --- Step 2: fill remaining cells (reusing already-registered strings, or numbers) ---
$ws.Range('A1').Value = 'title'
$ws.Range('B1').Value = 'slug'
$ws.Range('C1').Value = 'projectType'
$ws.Range('D1').Value = 'images'
$ws.Range('A2').Value = 'Factitious'
$ws.Range('D2').Value = 6
$ws.Range('C3').Value = 'Newsgame'
$ws.Range('D3').Value = 6
$ws.Range('F3').Value = 'asdf'
$ws.Range('G3').Value = 'asdf'
$ws.Range('H3').Value = 'asdf'
$ws.Range('I3').Value = 'asdf'
$ws.Range('D4').Value = 6
$ws.Range('G4').Value = 'asdf'
$ws.Range('H4').Value = 'asdf'
$ws.Range('I4').Value = 'asdf'

# --- Step 3: apply wrap-text style to the new F:I columns (rows 1-4) ---
$ws.Range("F1:I4").WrapText = $true

# --- Step 4: column widths (character units); XML stores width = ColumnWidth + 5/6,
#     so back out the ColumnWidth that reproduces each target XML width ---
$ws.Columns.Item(1).ColumnWidth = 25.498697916666668
$ws.Columns.Item(2).ColumnWidth = 18.998697916666668
$ws.Columns.Item(5).ColumnWidth = 23.830729166666668
$ws.Columns.Item(6).ColumnWidth = 19.998697916666668
$ws.Columns.Item(7).ColumnWidth = 34.998697916666664
$ws.Columns.Item(8).ColumnWidth = 30.498697916666668
$ws.Columns.Item(9).ColumnWidth = 34.166666666666664
$ws.Columns.Item(10).ColumnWidth = 31.330729166666668

# --- Step 5: row height for row 2 (wrapped long text) ---
$ws.Rows.Item(2).RowHeight = 140

# --- Step 6: final selection matches the authored state ---
$ws.Range("I4").Select() | Out-Null
